$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rule R30's "From" threshold (cell C10) changes from 18 to 100
$ws.Range("C10").Value = 100
